$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of old Taxonsorteringsordning (column B) values to new ones
$map = @{
    79239 = 79243
    91824 = 91828
    57880 = 57884
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($null -ne $val) {
        $key = [int]$val
        if ($map.ContainsKey($key)) {
            $cell.Value2 = $map[$key]
        }
    }
}
